# This workbook is a pure "results" data dump (no formulas) produced by an
# external simulation run. The commit re-uploads the workbook after a
# recalculation of the underlying model, which changed a handful of data
# cells (and added a few previously-empty/zero cells) across several sheets.
# We simply poke the new literal values into the corresponding cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "pcroprep" (sheet1.xml)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("pcroprep")

$ws1.Range("D35").Value = 1.4105556681109279
$ws1.Range("F35").Value = 5.1292933385851915
$ws1.Range("G35").Value = -236.27070666141481

$ws1.Range("D39").Value = 1177.392608731664
$ws1.Range("F39").Value = 371.00188407754888
$ws1.Range("G39").Value = -435.59811592245109

# ---------------------------------------------------------------------
# Sheet "pdietrep" (sheet4.xml)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("pdietrep")

$ws4.Range("E6").Value = 1502.8167900564476
$ws4.Range("F6").Value = -654.96872414923519
$ws4.Range("G6").Value = 69.646254466105248

$ws4.Range("E7").Value = 57.338852480156461
$ws4.Range("F7").Value = -12.789338396194097
$ws4.Range("G7").Value = 81.762914120023225

$ws4.Range("E8").Value = 23.8967783091258
$ws4.Range("F8").Value = -40.836787117044672
$ws4.Range("G8").Value = 36.915591087564003

$ws4.Range("E9").Value = 247.79858836028561
$ws4.Range("F9").Value = -75.869238770566739
$ws4.Range("G9").Value = 76.559536533764188

# ---------------------------------------------------------------------
# Sheet "pradar" (sheet5.xml)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("pradar")

$ws5.Range("D15").Value = 5.1292933385851915
$ws5.Range("E15").Value = 2.1248108279143296
$ws5.Range("F15").Value = -236.27070666141481

# ---------------------------------------------------------------------
# Sheet "plandrep" (sheet6.xml)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("plandrep")

$ws6.Range("S11").Value = 10.434000000000033

# ---------------------------------------------------------------------
# Sheet "plaborrep" (sheet7.xml)
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("plaborrep")

$ws7.Range("R3").Value = 0.0085037100000000098
$ws7.Range("AF3").Value = 1.3673698530208134

# ---------------------------------------------------------------------
# Sheet "pfertrep" (sheet8.xml)
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("pfertrep")

$ws8.Range("S5").Value = 2869.350000000009
$ws8.Range("Z5").Value = 393192.13386399997

$ws8.Range("S6").Value = 1721.6100000000054
$ws8.Range("Z6").Value = 493017.38162400003

$ws8.Range("S7").Value = 2744.1420000000085
$ws8.Range("Z7").Value = 458503.66967199993
